$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.932.79'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '2.354.12'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.65'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -2.49%  '
$ws.Range("D9").Value = '2.366.81'
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0958'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.28%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.78'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.321'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("D14").Value = '2.774.46'
$ws.Range("E14").Value = '  -1.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.72'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '55.897.21'
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000130'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").Value = '2.384.85'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.89'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.65%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.03'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '308.78'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.48%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.24'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.40'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.372'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("E27").Value = '  -3.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.22'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.57'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.84%  '
$ws.Range("D30").Value = '0.0₃0711'
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.78'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.08'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.68'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.18'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -1.52%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.864'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.70'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.27'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("E41").Value = '  -2.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.374'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.93'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '125.33'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.555'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0894'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '242.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0481'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.96'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  -1.24%  '
